$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.409.05"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "2.157.57"
$ws.Range("E3").Value = "  +3.22%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.53"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.60%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0859"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.29%  "

$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.62%  "

$ws.Range("D13").Value = "2.478.82"
$ws.Range("E13").Value = "  +3.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.13%  "

$ws.Range("E15").Value = "  +2.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("D17").Value = "2.147.36"
$ws.Range("E17").Value = "  +2.40%  "

$ws.Range("D18").Value = "39.519.25"
$ws.Range("E18").Value = "  +2.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.78%  "

$ws.Range("D21").Value = "0.0₃0854"
$ws.Range("E21").Value = "  +1.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.33%  "

$ws.Range("E30").Value = "  -1.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.42%  "

$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.48%  "

$ws.Range("E41").Value = "  +3.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.18%  "

$ws.Range("D43").Value = "1.535.07"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("E44").Value = "  +5.95%  "

$ws.Range("E45").Value = "  +7.39%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0923"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("E49").Value = "  +1.87%  "

$ws.Range("D50").Value = "2.361.66"
$ws.Range("E50").Value = "  +3.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
